$d = $word.ActiveDocument

# Merge the title runs "AutoCam" + " Interrupts" into a single run of text,
# which also removes the proofErr spell-check markers surrounding "AutoCam".
$d.Content.Find.Execute("AutoCam Interrupts", $false, $false, $false, $false, $false,
                         $true, 1, $false, "AutoCam Interrupts", 2) | Out-Null

# Fill in the first data row of the table.
$table = $d.Tables.Item(1)
$table.Cell(2, 1).Range.Text = "TimerCompare1"
$table.Cell(2, 2).Range.Text = "For the camera shooting activities, counts in seconds also if initiated"
